# Weekly update for "Hortaliza, Vega Central Mapocho de Santiago - Jengibre":
# a new weekly price record is inserted as the new top data row (row 53),
# pushing the previously-existing rows 53:64 down to 54:65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 53 (shifts 53:64 -> 54:65, grows the used range to R65)
$ws.Rows.Item(53).Insert()

# Populate the new row with the latest weekly observation
$ws.Cells.Item(53, 1).Value = 9
$ws.Cells.Item(53, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(53, 3).Value = 'Metropolitana'
$ws.Cells.Item(53, 4).Value = 44508
$ws.Cells.Item(53, 5).Value = 13
$ws.Cells.Item(53, 6).Value = 100114007
$ws.Cells.Item(53, 7).Value = 'Jengibre'
$ws.Cells.Item(53, 8).Value = 'Sin especificar'
$ws.Cells.Item(53, 9).Value = 'Primera'
$ws.Cells.Item(53, 10).Value = 1150
$ws.Cells.Item(53, 11).Value = 15000
$ws.Cells.Item(53, 12).Value = 16000
$ws.Cells.Item(53, 13).Value = 15500
$ws.Cells.Item(53, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(53, 15).Value = 'Perú'
$ws.Cells.Item(53, 16).Value = 1192
$ws.Cells.Item(53, 17).Value = 13
$ws.Cells.Item(53, 18).Value = 'Hortaliza'
